$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2683.3333
$ws.Range("H132").Value = 6178822
$ws.Range("I132").Value = 8551173
$ws.Range("J132").Value = 10709.533
$ws.Range("K132").Value = 25653519
$ws.Range("L132").Value = 32128.599
$ws.Range("M132").Value = -25650989
$ws.Range("N132").Value = -37188.599
$ws.Range("H137").Value = 1266.7097
$ws.Range("I137").Value = 767
$ws.Range("J137").Value = 1873.5
$ws.Range("K137").Value = 2301
$ws.Range("L137").Value = 5620.5
$ws.Range("M137").Value = 249
$ws.Range("N137").Value = -10720.5
$ws.Range("H138").Value = 1584.3636
$ws.Range("I138").Value = 1105.1428
$ws.Range("J138").Value = 1713.3846
$ws.Range("K138").Value = 3315.4284
$ws.Range("L138").Value = 5140.1538
$ws.Range("M138").Value = 1824.5716
$ws.Range("N138").Value = -15420.1538

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 29756.25
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 29756.25
$ws.Range("K23").Value = 0
$ws.Range("L23").ClearContents()
$ws.Range("M23").Value = 29756.25
$ws.Range("N23").Value = -30274.25
$ws.Range("H32").Value = 6858.788
$ws.Range("I32").Value = 6339.3447
$ws.Range("K32").Value = 6339.3447
$ws.Range("M32").Value = -6052.3447
$ws.Range("H74").Value = 1609.9231
$ws.Range("I74").Value = 859.6
$ws.Range("J74").Value = 2633.0908
$ws.Range("K74").Value = 859.6
$ws.Range("L74").Value = 2633.0908
$ws.Range("M74").Value = 14.39999999999998
$ws.Range("N74").Value = -4381.0908
$ws.Range("H77").Value = 1609.9231
$ws.Range("I77").Value = 859.6
$ws.Range("J77").Value = 2633.0908
$ws.Range("K77").Value = 4298
$ws.Range("L77").Value = 13165.454
$ws.Range("M77").Value = 70
$ws.Range("N77").Value = -21901.454
$ws.Range("H132").Value = 1310.2034
$ws.Range("I132").Value = 1048.4048
$ws.Range("K132").Value = 3145.2144
$ws.Range("M132").Value = -615.2143999999998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 31250990
$ws.Range("I99").Value = 35715204
$ws.Range("K99").Value = 35715204
$ws.Range("M99").Value = -35713706

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1359.6666
$ws.Range("I31").Value = 1171.6666
$ws.Range("J31").Value = 1735.6666
$ws.Range("K31").Value = 1171.6666
$ws.Range("L31").Value = 1735.6666
$ws.Range("M31").Value = -876.6666
$ws.Range("N31").Value = -2325.6666
$ws.Range("H34").Value = 1359.6666
$ws.Range("I34").Value = 1171.6666
$ws.Range("J34").Value = 1735.6666
$ws.Range("K34").Value = 1171.6666
$ws.Range("L34").Value = 1735.6666
$ws.Range("M34").Value = -969.6666
$ws.Range("N34").Value = -2139.6666
$ws.Range("H132").Value = 1594.3405
$ws.Range("I132").Value = 1216.8948
$ws.Range("J132").Value = 3188
$ws.Range("K132").Value = 3650.6844
$ws.Range("L132").Value = 9564
$ws.Range("M132").Value = -1120.6844
$ws.Range("N132").Value = -14624
$ws.Range("H134").Value = 679.79486
$ws.Range("I134").Value = 587.9375
$ws.Range("J134").Value = 1099.7142
$ws.Range("K134").Value = 1763.8125
$ws.Range("L134").Value = 3299.1426
$ws.Range("M134").Value = 771.1875
$ws.Range("N134").Value = -8369.142599999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 73865
$ws.Range("I109").Value = 167380
$ws.Range("K109").Value = 502140
$ws.Range("M109").Value = -501100

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2171.6667
$ws.Range("I82").Value = 2500
$ws.Range("J82").Value = 1952.7778
$ws.Range("K82").Value = 2500
$ws.Range("L82").Value = 1952.7778
$ws.Range("M82").Value = -2139
$ws.Range("N82").Value = -2674.7778
$ws.Range("H85").Value = 2171.6667
$ws.Range("I85").Value = 2500
$ws.Range("J85").Value = 1952.7778
$ws.Range("K85").Value = 2500
$ws.Range("L85").Value = 1952.7778
$ws.Range("M85").Value = -1252
$ws.Range("N85").Value = -4448.7778
$ws.Range("H124").ClearContents()
$ws.Range("I124").ClearContents()
$ws.Range("J124").ClearContents()
$ws.Range("K124").ClearContents()
$ws.Range("L124").ClearContents()
$ws.Range("H125").ClearContents()
$ws.Range("I125").ClearContents()
$ws.Range("J125").ClearContents()
$ws.Range("K125").ClearContents()
$ws.Range("L125").ClearContents()
$ws.Range("N125").ClearContents()
$ws.Range("H127").ClearContents()
$ws.Range("I127").ClearContents()
$ws.Range("J127").ClearContents()
$ws.Range("K127").ClearContents()
$ws.Range("L127").ClearContents()
$ws.Range("N127").ClearContents()
$ws.Range("H128").ClearContents()
$ws.Range("I128").ClearContents()
$ws.Range("J128").ClearContents()
$ws.Range("K128").ClearContents()
$ws.Range("L128").ClearContents()
$ws.Range("N128").ClearContents()
$ws.Range("H129").ClearContents()
$ws.Range("I129").ClearContents()
$ws.Range("J129").ClearContents()
$ws.Range("K129").ClearContents()
$ws.Range("L129").ClearContents()
$ws.Range("H130").ClearContents()
$ws.Range("I130").ClearContents()
$ws.Range("J130").ClearContents()
$ws.Range("K130").ClearContents()
$ws.Range("L130").ClearContents()
$ws.Range("N130").ClearContents()
$ws.Range("H131").ClearContents()
$ws.Range("I131").ClearContents()
$ws.Range("J131").ClearContents()
$ws.Range("K131").ClearContents()
$ws.Range("L131").ClearContents()
$ws.Range("N131").ClearContents()
$ws.Range("H132").ClearContents()
$ws.Range("I132").ClearContents()
$ws.Range("J132").ClearContents()
$ws.Range("K132").ClearContents()
$ws.Range("L132").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H133").ClearContents()
$ws.Range("I133").ClearContents()
$ws.Range("J133").ClearContents()
$ws.Range("K133").ClearContents()
$ws.Range("L133").ClearContents()
$ws.Range("N133").ClearContents()
$ws.Range("H134").ClearContents()
$ws.Range("I134").ClearContents()
$ws.Range("J134").ClearContents()
$ws.Range("K134").ClearContents()
$ws.Range("L134").ClearContents()
$ws.Range("N134").ClearContents()
$ws.Range("H135").ClearContents()
$ws.Range("I135").ClearContents()
$ws.Range("J135").ClearContents()
$ws.Range("K135").ClearContents()
$ws.Range("L135").ClearContents()
$ws.Range("N135").ClearContents()
$ws.Range("H136").ClearContents()
$ws.Range("I136").ClearContents()
$ws.Range("J136").ClearContents()
$ws.Range("K136").ClearContents()
$ws.Range("L136").ClearContents()
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()
$ws.Range("H137").ClearContents()
$ws.Range("I137").ClearContents()
$ws.Range("J137").ClearContents()
$ws.Range("K137").ClearContents()
$ws.Range("L137").ClearContents()
$ws.Range("N137").ClearContents()
$ws.Range("H138").ClearContents()
$ws.Range("I138").ClearContents()
$ws.Range("J138").ClearContents()
$ws.Range("K138").ClearContents()
$ws.Range("L138").ClearContents()
$ws.Range("H139").ClearContents()
$ws.Range("I139").ClearContents()
$ws.Range("J139").ClearContents()
$ws.Range("K139").ClearContents()
$ws.Range("L139").ClearContents()
$ws.Range("N139").ClearContents()
$ws.Range("H140").ClearContents()
$ws.Range("I140").ClearContents()
$ws.Range("J140").ClearContents()
$ws.Range("K140").ClearContents()
$ws.Range("L140").ClearContents()
$ws.Range("N140").ClearContents()
$ws.Range("H141").ClearContents()
$ws.Range("I141").ClearContents()
$ws.Range("J141").ClearContents()
$ws.Range("K141").ClearContents()
$ws.Range("L141").ClearContents()
$ws.Range("N141").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 35530
$ws.Range("J75").Value = 35530
$ws.Range("L75").Value = 35530
$ws.Range("N75").Value = -37402
$ws.Range("H78").Value = 35530
$ws.Range("J78").Value = 35530
$ws.Range("L78").Value = 106590
$ws.Range("N78").Value = -115950
$ws.Range("H132").Value = 2218.4707
$ws.Range("I132").Value = 2157.8333
$ws.Range("J132").Value = 2673.25
$ws.Range("K132").Value = 6473.499899999999
$ws.Range("L132").Value = 8019.75
$ws.Range("M132").Value = -3943.499899999999
$ws.Range("N132").Value = -13079.75
$ws.Range("H136").Value = 574.4375
$ws.Range("I136").Value = 253.07692
$ws.Range("K136").Value = 759.23076
$ws.Range("M136").Value = 1790.76924
